$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 6474
    5  = 42
    6  = 1944
    7  = 1477
    9  = 992
    10 = 345
    11 = 2
    12 = 5614
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
